$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "as of" serial number in I1 (e.g. a scan/report date counter)
$ws.Range("I1").Value = 33557

# Clear out the old "X" mark and bill amount on Blue Cross Blue Shield of Illinois (row 3)
$ws.Range("B3:D3").ClearContents()

# Clear out the old "X" mark, invoice number and bill amount on Hernandez Lawn Service (row 12)
$ws.Range("B12:D12").ClearContents()

# Mark CNA Insurance (row 6) as the vendor to use, with its bill amount
$ws.Range("B6").Value = "x"
$ws.Range("D6").Value = 87836.9

# Update the Nicor Gas bill total formula (row 16)
$ws.Range("I16").Formula = "=156.05+175.85+262.61"

# Move the active selection to I1, matching where the scan tool left the cursor
$ws.Range("I1").Select()
